$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new expense entry added (icon URL, amount, date)
$ws.Range("A2").Value = "https://cdn.jsdelivr.net/npm/emoji-datasource-apple/img/apple/64/1f386.png"
$ws.Range("C2").Value = 1233
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "5/30/2025"

# Row 3: previous row 2 values
$ws.Range("A3").Value = "icons2"
$ws.Range("C3").Value = 5000
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4/23/2025"

# Row 4: previous row 3 values
$ws.Range("A4").Value = "https://cdn.jsdelivr.net/npm/emoji-datasource-apple/img/apple/64/1f977.png"
$ws.Range("C4").Value = 3331
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4/15/2025"
